$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.077.37"
$ws.Range("E2").Value = "  -4.13%  "

$ws.Range("D3").Value = "2.895.73"
$ws.Range("E3").Value = "  -3.56%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.38%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("D9").Value = "2.897.22"
$ws.Range("E9").Value = "  -3.41%  "

$ws.Range("E10").Value = "  -4.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.60%  "

$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("D13").Value = "3.402.65"
$ws.Range("E13").Value = "  -3.34%  "

$ws.Range("E14").Value = "  +1.13%  "

$ws.Range("D15").Value = "60.225.98"
$ws.Range("E15").Value = "  -4.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.18%  "

$ws.Range("D17").Value = "2.897.48"
$ws.Range("E17").Value = "  -3.70%  "

$ws.Range("E18").Value = "  -6.12%  "

$ws.Range("E19").Value = "  -2.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.79%  "

$ws.Range("D26").Value = "3.013.47"
$ws.Range("E26").Value = "  -4.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.447"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.174"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.95%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.55%  "

$ws.Range("D31").Value = "0.0₃0850"
$ws.Range("E31").Value = "  -11.35%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  -4.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.21%  "

$ws.Range("E37").Value = "  -7.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.983"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.09%  "

$ws.Range("E39").Value = "  -7.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").Value = "2.333.78"
$ws.Range("E41").Value = "  -5.71%  "

$ws.Range("E42").Value = "  -6.67%  "

$ws.Range("E43").Value = "  -5.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.643"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0567"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.40%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  -2.89%  "

$ws.Range("E49").Value = "  -1.36%  "

$ws.Range("E50").Value = "  -5.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0926"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.67%  "
